$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15; this shifts the existing rows 15..130 down to 16..131
# and inherits cell formatting (e.g. the date number format on column D) from the
# row above, same as Excel's native Insert behaviour.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with its data (a new weekly price observation).
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44490
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 100112010
$ws.Range("G15").Value = "Achicoria"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 5500
$ws.Range("L15").Value = 6000
$ws.Range("M15").Value = 5750
$ws.Range("N15").Value = "$/caja 16 unidades"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 359
$ws.Range("Q15").Value = 16
$ws.Range("R15").Value = "Hortaliza"
